$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the schedule row for Joshua Aguirre (row 4): B-F = 5pm-MN, G-H = 8am-MN
$ws.Range("B4").Value = "5pm-MN"
$ws.Range("C4").Value = "5pm-MN"
$ws.Range("D4").Value = "5pm-MN"
$ws.Range("E4").Value = "5pm-MN"
$ws.Range("F4").Value = "5pm-MN"
$ws.Range("G4").Value = "8am-MN"
$ws.Range("H4").Value = "8am-MN"

# Update the active selection to match where the edit ended up
$null = $ws.Range("G4").Select()
